$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 1213.6765
$ws.Range("J17").Value = 1213.6765
$ws.Range("L17").Value = 3641.0295
$ws.Range("N17").Value = -3977.0295

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 1346.1666
$ws.Range("I18").Value = 1333.3334
$ws.Range("J18").Value = 1384.6666
$ws.Range("K18").Value = 1333.3334
$ws.Range("L18").Value = 1384.6666
$ws.Range("M18").Value = -1049.3334
$ws.Range("N18").Value = -1952.6666

# Row 37 (Leve Item ID 4621)
$ws.Range("H37").Value = 900
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 2427.7144
$ws.Range("I51").Value = 1798.5
$ws.Range("J51").Value = 3266.6667
$ws.Range("K51").Value = 1798.5
$ws.Range("L51").Value = 3266.6667
$ws.Range("M51").Value = -1314.5
$ws.Range("N51").Value = -4234.6667

# Row 54 (Leve Item ID 2174)
$ws.Range("H54").Value = 500
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 7425.5415
$ws.Range("I98").Value = 8018.409
$ws.Range("K98").Value = 8018.409
$ws.Range("M98").Value = -6520.409

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 5962.5
$ws.Range("I111").Value = 5962.5
$ws.Range("K111").Value = 17887.5
$ws.Range("M111").Value = -14820.5

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 7425.5415
$ws.Range("I122").Value = 8018.409
$ws.Range("K122").Value = 24055.227
$ws.Range("M122").Value = -21605.227

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 849.6391599999999
$ws.Range("J129").Value = 918.97675
$ws.Range("L129").Value = 2756.93025
$ws.Range("N129").Value = -12756.93025

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3553.026
$ws.Range("I32").Value = 3243.6572
$ws.Range("K32").Value = 3243.6572
$ws.Range("M32").Value = -2956.6572

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1843.2572
$ws.Range("I132").Value = 1758.3334
$ws.Range("J132").Value = 1933.1765
$ws.Range("K132").Value = 5275.0002
$ws.Range("L132").Value = 5799.529500000001
$ws.Range("M132").Value = -2745.0002
$ws.Range("N132").Value = -10859.5295

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1809.3334
$ws.Range("I20").Value = 1795.4546
$ws.Range("K20").Value = 1795.4546
$ws.Range("M20").Value = -1548.4546

# Row 116 (Leve Item ID 26113)
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 5894.1377
$ws.Range("I134").Value = 1690.85
$ws.Range("K134").Value = 5072.549999999999
$ws.Range("M134").Value = -2537.549999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1566.3334
$ws.Range("I31").Value = 1499.5
$ws.Range("J31").Value = 1700
$ws.Range("K31").Value = 1499.5
$ws.Range("L31").Value = 1700
$ws.Range("M31").Value = -1204.5
$ws.Range("N31").Value = -2290

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1566.3334
$ws.Range("I34").Value = 1499.5
$ws.Range("J34").Value = 1700
$ws.Range("K34").Value = 1499.5
$ws.Range("L34").Value = 1700
$ws.Range("M34").Value = -1297.5
$ws.Range("N34").Value = -2104

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 3756
$ws.Range("I58").Value = 3756
$ws.Range("K58").Value = 3756
$ws.Range("M58").Value = -3553

# Row 100 (Leve Item ID 34388)
$ws.Range("H100").Value = 61633.332
$ws.Range("J100").Value = 61633.332
$ws.Range("L100").Value = 61633.332
$ws.Range("N100").Value = -63797.332

# Row 114 (Leve Item ID 27112)
$ws.Range("H114").Value = 24796
$ws.Range("J114").Value = 25745
$ws.Range("L114").Value = 25745
$ws.Range("N114").Value = -34423

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 948.1429000000001
$ws.Range("I134").Value = 786.7857
$ws.Range("K134").Value = 2360.3571
$ws.Range("M134").Value = 174.6428999999998

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 3756
$ws.Range("I136").Value = 3756
$ws.Range("K136").Value = 11268
$ws.Range("M136").Value = -8718

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1416
$ws.Range("I5").Value = 1684.762
$ws.Range("J5").Value = 788.8889
$ws.Range("K5").Value = 5054.286
$ws.Range("L5").Value = 2366.6667
$ws.Range("M5").Value = -4942.286
$ws.Range("N5").Value = -2590.6667

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1416
$ws.Range("I135").Value = 1684.762
$ws.Range("J135").Value = 788.8889
$ws.Range("K135").Value = 15162.858
$ws.Range("L135").Value = 7100.0001
$ws.Range("M135").Value = -12627.858
$ws.Range("N135").Value = -12170.0001

# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 33205.11
$ws.Range("J140").Value = 2835.2273
$ws.Range("L140").Value = 8505.6819
$ws.Range("N140").Value = -18865.6819

$ws = $wb.Worksheets.Item("GSM")
# Row 62 (Leve Item ID 11983)
$ws.Range("H62").Value = 20085
$ws.Range("J62").Value = 20085
$ws.Range("L62").Value = 20085
$ws.Range("N62").Value = -21457

# Row 65 (Leve Item ID 11983)
$ws.Range("H65").Value = 20085
$ws.Range("J65").Value = 20085
$ws.Range("L65").Value = 60255
$ws.Range("N65").Value = -67119

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 30009828
$ws.Range("I70").Value = 27791996
$ws.Range("J70").Value = 33336576
$ws.Range("K70").Value = 27791996
$ws.Range("L70").Value = 33336576
$ws.Range("M70").Value = -27791726
$ws.Range("N70").Value = -33337116

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 30009828
$ws.Range("I73").Value = 27791996
$ws.Range("J73").Value = 33336576
$ws.Range("K73").Value = 27791996
$ws.Range("L73").Value = 33336576
$ws.Range("M73").Value = -27791060
$ws.Range("N73").Value = -33338448

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 3730.647
$ws.Range("I102").Value = 2515.1738
$ws.Range("J102").Value = 6272.091
$ws.Range("K102").Value = 2515.1738
$ws.Range("L102").Value = 6272.091
$ws.Range("M102").Value = -893.1738
$ws.Range("N102").Value = -9516.091

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3042.7144
$ws.Range("I122").Value = 2716.6667
$ws.Range("K122").Value = 8150.000100000001
$ws.Range("M122").Value = -5700.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 5914
$ws.Range("I40").Value = 3404
$ws.Range("J40").Value = 6332.3335
$ws.Range("K40").Value = 3404
$ws.Range("L40").Value = 6332.3335
$ws.Range("M40").Value = -3268
$ws.Range("N40").Value = -6604.3335

# Row 63 (Leve Item ID 12006)
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 20000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -21498

# Row 64 (Leve Item ID 10810)
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 66 (Leve Item ID 12006)
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 60000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -67488

# Row 67 (Leve Item ID 10810)
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1870.381
$ws.Range("I132").Value = 1705
$ws.Range("K132").Value = 5115
$ws.Range("M132").Value = -2585

Write-Output "Applied Kujata_Profits market price updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
